# Atualização automática de pedidos - 30/05/2025 07:32
# Adds the new request REQ-006 to both the "Pedidos" and "Itens" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Pedidos": append order header row (row 7)
# ---------------------------------------------------------------------
$pedidos = $wb.Worksheets.Item("Pedidos")

$pedidos.Range("A7").Value = "REQ-006"
$pedidos.Range("B7").Value = "30/05/2025 07:32"
$pedidos.Range("C7").Value = "Renault"

# RACK is stored as text ("1") for this row, unlike the numeric values
# used above it, so force a text number format before assigning it.
$pedidos.Range("D7").NumberFormat = "@"
$pedidos.Range("D7").Value = "1"
$pedidos.Range("D7").Style = "Normal"

$pedidos.Range("E7").Value = "R01-LA-B2"
$pedidos.Range("F7").Value = "washington vieira"
$pedidos.Range("H7").Value = "Pendente"

# ---------------------------------------------------------------------
# Sheet "Itens": append item detail row (row 7)
# ---------------------------------------------------------------------
$itens = $wb.Worksheets.Item("Itens")

$itens.Range("A7").Value = "REQ-006"
$itens.Range("B7").Value = "A3ZPA-1.0-GY"
$itens.Range("C7").Value = "180DN106041"

# seccao is stored as text ("1.0") for this row, unlike the numeric
# values used above it, so force a text number format before assigning it.
$itens.Range("D7").NumberFormat = "@"
$itens.Range("D7").Value = "1.0"
$itens.Range("D7").Style = "Normal"

$itens.Range("E7").Value = "GY"
$itens.Range("F7").Value = 1
